# Insert a new weekly price record as row 147, pushing the existing
# rows 147-235 down to 148-236 (dimension grows from A1:T235 to A1:T236).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 147, shifting rows 147:235 down.
$ws.Rows("147:147").Insert()

# Populate the newly inserted row 147 with the new weekly record.
$ws.Range("A147").Value = 4
$ws.Range("B147").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C147").Value = "Los Lagos"
$ws.Range("D147").Value = 44438
$ws.Range("E147").Value = 10
$ws.Range("F147").Value = "Fruta"
$ws.Range("G147").Value = 100108
$ws.Range("H147").Value = "Tropicales y subtropicales"
$ws.Range("I147").Value = 100108006
$ws.Range("J147").Value = "Plátano"
$ws.Range("K147").Value = "Sin especificar"
$ws.Range("L147").Value = "Primera Pintón"
$ws.Range("M147").Value = 300
$ws.Range("N147").Value = 21000
$ws.Range("O147").Value = 21000
$ws.Range("P147").Value = 21000
$ws.Range("Q147").Value = "$/caja 20 kilos"
$ws.Range("R147").Value = "Ecuador"
$ws.Range("S147").Value = 1050
$ws.Range("T147").Value = 20
